$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 306.75
$ws.Range("I38").Value = 70
$ws.Range("J38").Value = 475.85715
$ws.Range("K38").Value = 210
$ws.Range("L38").Value = 1427.57145
$ws.Range("M38").Value = 162
$ws.Range("N38").Value = -2171.57145
$ws.Range("H98").Value = 1023056.06
$ws.Range("I98").Value = 1023056.06
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1023056.06
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H121").Value = 861
$ws.Range("J121").Value = 1201.6666
$ws.Range("L121").Value = 3604.9998
$ws.Range("N121").Value = -7098.9998
$ws.Range("H122").Value = 1023056.06
$ws.Range("I122").Value = 1023056.06
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3069168.18
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H138").Value = 11823376
$ws.Range("I138").Value = 2648441.8
$ws.Range("J138").Value = 19233900
$ws.Range("K138").Value = 7945325.399999999
$ws.Range("L138").Value = 57701700
$ws.Range("M138").Value = -7940185.399999999
$ws.Range("N138").Value = -57711980
$ws.Range("H141").Value = 4102.684
$ws.Range("I141").Value = 2327.889
$ws.Range("J141").Value = 8459
$ws.Range("K141").Value = 6983.667
$ws.Range("L141").Value = 25377
$ws.Range("M141").Value = -1803.667
$ws.Range("N141").Value = -35737

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4457.913
$ws.Range("I122").Value = 4374.364
$ws.Range("J122").Value = 4534.5
$ws.Range("K122").Value = 13123.092
$ws.Range("L122").Value = 13603.5
$ws.Range("M122").Value = -10673.092
$ws.Range("N122").Value = -18503.5
$ws.Range("H132").Value = 2339.0417
$ws.Range("I132").Value = 1949.2683
$ws.Range("J132").Value = 4622
$ws.Range("K132").Value = 5847.8049
$ws.Range("L132").Value = 13866
$ws.Range("M132").Value = -3317.8049
$ws.Range("N132").Value = -18926

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1519.125
$ws.Range("I20").Value = 1311.7778
$ws.Range("J20").Value = 1785.7142
$ws.Range("K20").Value = 1311.7778
$ws.Range("L20").Value = 1785.7142
$ws.Range("M20").Value = -1064.7778
$ws.Range("N20").Value = -2279.7142
$ws.Range("H29").Value = 665.3333
$ws.Range("I29").Value = 665.3333
$ws.Range("K29").Value = 665.3333
$ws.Range("M29").Value = -376.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5778.604
$ws.Range("I31").Value = 1960.762
$ws.Range("J31").Value = 8748.037
$ws.Range("K31").Value = 1960.762
$ws.Range("L31").Value = 8748.037
$ws.Range("M31").Value = -1665.762
$ws.Range("N31").Value = -9338.037
$ws.Range("H34").Value = 5778.604
$ws.Range("I34").Value = 1960.762
$ws.Range("J34").Value = 8748.037
$ws.Range("K34").Value = 1960.762
$ws.Range("L34").Value = 8748.037
$ws.Range("M34").Value = -1758.762
$ws.Range("N34").Value = -9152.037
$ws.Range("H58").Value = 2199.6956
$ws.Range("I58").Value = 1293.875
$ws.Range("J58").Value = 4270.143
$ws.Range("K58").Value = 1293.875
$ws.Range("L58").Value = 4270.143
$ws.Range("M58").Value = -1090.875
$ws.Range("N58").Value = -4676.143
$ws.Range("H122").Value = 1401.875
$ws.Range("I122").Value = 1235.8334
$ws.Range("K122").Value = 3707.5002
$ws.Range("M122").Value = -1257.5002
$ws.Range("H136").Value = 2199.6956
$ws.Range("I136").Value = 1293.875
$ws.Range("J136").Value = 4270.143
$ws.Range("K136").Value = 3881.625
$ws.Range("L136").Value = 12810.429
$ws.Range("M136").Value = -1331.625
$ws.Range("N136").Value = -17910.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 141.26086
$ws.Range("I4").Value = 86.59999999999999
$ws.Range("J4").Value = 243.75
$ws.Range("K4").Value = 259.8
$ws.Range("L4").Value = 731.25
$ws.Range("M4").Value = -147.8
$ws.Range("N4").Value = -955.25
$ws.Range("H6").Value = 213.55556
$ws.Range("I6").Value = 103.666664
$ws.Range("J6").Value = 433.33334
$ws.Range("K6").Value = 310.999992
$ws.Range("L6").Value = 1300.00002
$ws.Range("M6").Value = -197.999992
$ws.Range("N6").Value = -1526.00002
$ws.Range("H7").Value = 80.333336
$ws.Range("I7").Value = 80.333336
$ws.Range("K7").Value = 241.000008
$ws.Range("M7").Value = -129.000008
$ws.Range("H29").Value = 2128.6
$ws.Range("J29").Value = 160.75
$ws.Range("L29").Value = 482.25
$ws.Range("N29").Value = -1036.25
$ws.Range("H113").Value = 33333864
$ws.Range("J113").Value = 33333864
$ws.Range("L113").Value = 100001592
$ws.Range("N113").Value = -100005932
$ws.Range("H122").Value = 2215.889
$ws.Range("I122").Value = 304.5
$ws.Range("J122").Value = 2762
$ws.Range("K122").Value = 2740.5
$ws.Range("L122").Value = 24858
$ws.Range("M122").Value = -290.5
$ws.Range("N122").Value = -29758
$ws.Range("H137").Value = 5319506.5
$ws.Range("I137").Value = 12503320
$ws.Range("J137").Value = 94915
$ws.Range("K137").Value = 37509960
$ws.Range("L137").Value = 284745
$ws.Range("M137").Value = -37504860
$ws.Range("N137").Value = -294945

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10882.5
$ws.Range("I70").Value = 34900
$ws.Range("J70").Value = 6079
$ws.Range("K70").Value = 34900
$ws.Range("L70").Value = 6079
$ws.Range("M70").Value = -34630
$ws.Range("N70").Value = -6619
$ws.Range("H73").Value = 10882.5
$ws.Range("I73").Value = 34900
$ws.Range("J73").Value = 6079
$ws.Range("K73").Value = 34900
$ws.Range("L73").Value = 6079
$ws.Range("M73").Value = -33964
$ws.Range("N73").Value = -7951
$ws.Range("H80").Value = 2739.6072
$ws.Range("I80").Value = 2662.3809
$ws.Range("J80").Value = 2971.2856
$ws.Range("K80").Value = 2662.3809
$ws.Range("L80").Value = 2971.2856
$ws.Range("M80").Value = -1664.3809
$ws.Range("N80").Value = -4967.2856
$ws.Range("H83").Value = 2739.6072
$ws.Range("I83").Value = 2662.3809
$ws.Range("J83").Value = 2971.2856
$ws.Range("K83").Value = 13311.9045
$ws.Range("L83").Value = 14856.428
$ws.Range("M83").Value = -8319.904500000001
$ws.Range("N83").Value = -24840.428
$ws.Range("H122").Value = 2370.5715
$ws.Range("I122").Value = 2437
$ws.Range("J122").Value = 2320.75
$ws.Range("K122").Value = 7311
$ws.Range("L122").Value = 6962.25
$ws.Range("M122").Value = -4861
$ws.Range("N122").Value = -11862.25
$ws.Range("H132").Value = 3758.353
$ws.Range("I132").Value = 3870.9092
$ws.Range("K132").Value = 11612.7276
$ws.Range("M132").Value = -9082.7276

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3045.8635
$ws.Range("I7").Value = 2099.8333
$ws.Range("J7").Value = 3400.625
$ws.Range("K7").Value = 2099.8333
$ws.Range("L7").Value = 3400.625
$ws.Range("M7").Value = -1987.8333
$ws.Range("N7").Value = -3624.625
$ws.Range("H98").Value = 19203.334
$ws.Range("J98").Value = 19203.334
$ws.Range("L98").Value = 19203.334
$ws.Range("N98").Value = -25193.334
$ws.Range("H126").Value = 3045.8635
$ws.Range("I126").Value = 2099.8333
$ws.Range("J126").Value = 3400.625
$ws.Range("K126").Value = 6299.499899999999
$ws.Range("L126").Value = 10201.875
$ws.Range("M126").Value = -3829.499899999999
$ws.Range("N126").Value = -15141.875
$ws.Range("H132").Value = 5267.727
$ws.Range("I132").Value = 3372
$ws.Range("J132").Value = 6351
$ws.Range("K132").Value = 10116
$ws.Range("L132").Value = 19053
$ws.Range("M132").Value = -7586
$ws.Range("N132").Value = -24113
$ws.Range("H133").Value = 52927.867
$ws.Range("J133").Value = 52927.867
$ws.Range("L133").Value = 52927.867
$ws.Range("N133").Value = -57987.867
$ws.Range("H136").Value = 2666.6511
$ws.Range("I136").Value = 1478.4872
$ws.Range("J136").Value = 14251.25
$ws.Range("K136").Value = 4435.461600000001
$ws.Range("L136").Value = 42753.75
$ws.Range("M136").Value = -1885.461600000001
$ws.Range("N136").Value = -47853.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 442.5
$ws.Range("I113").Value = 433.77777
$ws.Range("K113").Value = 1301.33331
$ws.Range("M113").Value = 868.66669
$ws.Range("H126").Value = 53595.367
$ws.Range("I126").Value = 83908.914
$ws.Range("J126").Value = 1629.2858
$ws.Range("K126").Value = 251726.742
$ws.Range("L126").Value = 4887.857400000001
$ws.Range("M126").Value = -249256.742
$ws.Range("N126").Value = -9827.857400000001
